$wb = $excel.ActiveWorkbook

# --- Rename "SmallCap" sheet to "NAV" ---
$navSheet = $wb.Worksheets.Item("SmallCap")
$navSheet.Name = "NAV"

# --- Update view/selection/zoom on the NAV sheet ---
$navSheet.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 120
$navSheet.Range("F13").Select() | Out-Null

# Column widths: B:C stay the same width but lose "best fit" sizing (now user fixed),
# and D:E split into two individually sized (but same-width) columns.
$navSheet.Columns.Item(2).ColumnWidth = 6.830729166666667
$navSheet.Columns.Item(3).ColumnWidth = 6.830729166666667
$navSheet.Columns.Item(4).ColumnWidth = 7.330729166666667
$navSheet.Columns.Item(5).ColumnWidth = 7.330729166666667

# --- Update view/selection on the 52Week sheet ---
$weekSheet = $wb.Worksheets.Item("52Week")
$weekSheet.Activate() | Out-Null
$weekSheet.Range("B4").Select() | Out-Null

# Give column A an explicit width (previously relying on the sheet default).
$weekSheet.Columns.Item(1).ColumnWidth = 12.998697916666666

# Re-activate the 52Week sheet/tab as the workbook's last active sheet.
$weekSheet.Activate() | Out-Null
